$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.692291333333333
$ws.Range("H2").Value = 8.076874
$ws.Range("I2").Value = 0.1057975873398775
$ws.Range("J2").Value = 0.1132142695545834
$ws.Range("M2").Value = 86.066935
$ws.Range("N2").Value = 258.200805
$ws.Range("O2").Value = 0.8916716774694496
$ws.Range("P2").Value = 0.8955157110805073
$ws.Range("Q2").Value = 231.7172631870633
$ws.Range("R2").Value = 2085.45536868357
$ws.Range("S2").Value = 0.09433671217556916
$ws.Range("T2").Value = 0.101385157104633
$ws.Range("G3").Value = 2.692291333333333
$ws.Range("H3").Value = 8.076874
$ws.Range("I3").Value = 0.1057975873398775
$ws.Range("J3").Value = 0.1132142695545834
$ws.Range("O3").Value = 0.003820894467605101
$ws.Range("P3").Value = 0.003837366502243974
$ws.Range("Q3").Value = 0.9929296077595555
$ws.Range("R3").Value = 8.936366469836001
$ws.Range("S3").Value = 0.0004042414161529054
$ws.Range("T3").Value = 0.0004344446455647781
$ws.Range("G4").Value = 2.692291333333333
$ws.Range("H4").Value = 8.076874
$ws.Range("I4").Value = 0.1057975873398775
$ws.Range("J4").Value = 0.1132142695545834
$ws.Range("M4").Value = 7.669867666666666
$ws.Range("N4").Value = 23.009603
$ws.Range("O4").Value = 0.07946145367329926
$ws.Range("P4").Value = 0.07980401529819077
$ws.Range("Q4").Value = 20.64951824678022
$ws.Range("R4").Value = 185.845664221022
$ws.Range("S4").Value = 0.008406830085154507
$ws.Range("T4").Value = 0.009034953299507468
$ws.Range("G5").Value = 2.692291333333333
$ws.Range("H5").Value = 8.076874
$ws.Range("I5").Value = 0.1057975873398775
$ws.Range("J5").Value = 0.1132142695545834
$ws.Range("M5").Value = 1.242987
$ws.Range("N5").Value = 2.485974
$ws.Range("O5").Value = 0.01287760861197995
$ws.Range("P5").Value = 0.008622083011467191
$ws.Range("Q5").Value = 3.346483127546
$ws.Range("R5").Value = 20.078898765276
$ws.Range("S5").Value = 0.001362419921854707
$ws.Range("T5").Value = 0.0009761428301822407
$ws.Range("G6").Value = 2.692291333333333
$ws.Range("H6").Value = 8.076874
$ws.Range("I6").Value = 0.1057975873398775
$ws.Range("J6").Value = 0.1132142695545834
$ws.Range("M6").Value = 1.174528666666667
$ws.Range("N6").Value = 3.523586
$ws.Range("O6").Value = 0.01216836577766621
$ws.Range("P6").Value = 0.01222082410759068
$ws.Range("Q6").Value = 3.162173350018222
$ws.Range("R6").Value = 28.459560150164
$ws.Range("S6").Value = 0.001287383741146218
$ws.Range("T6").Value = 0.001383571674695923
$ws.Range("I7").Value = 0.1778087899819636
$ws.Range("J7").Value = 0.1902736421911268
$ws.Range("M7").Value = 86.066935
$ws.Range("N7").Value = 258.200805
$ws.Range("O7").Value = 0.8916716774694496
$ws.Range("P7").Value = 0.8955157110805073
$ws.Range("Q7").Value = 389.4357822439134
$ws.Range("R7").Value = 3504.92204019522
$ws.Range("S7").Value = 0.1585470620320305
$ws.Range("T7").Value = 0.170393035986665
$ws.Range("I8").Value = 0.1778087899819636
$ws.Range("J8").Value = 0.1902736421911268
$ws.Range("O8").Value = 0.003820894467605101
$ws.Range("P8").Value = 0.003837366502243974
$ws.Range("S8").Value = 0.000679388621933642
$ws.Range("T8").Value = 0.0007301497008041858
$ws.Range("I9").Value = 0.1778087899819636
$ws.Range("J9").Value = 0.1902736421911268
$ws.Range("M9").Value = 7.669867666666666
$ws.Range("N9").Value = 23.009603
$ws.Range("O9").Value = 0.07946145367329926
$ws.Range("P9").Value = 0.07980401529819077
$ws.Range("Q9").Value = 34.70462744462355
$ws.Range("R9").Value = 312.341647001612
$ws.Range("S9").Value = 0.0141289449278572
$ws.Range("T9").Value = 0.01518460065226316
$ws.Range("I10").Value = 0.1778087899819636
$ws.Range("J10").Value = 0.1902736421911268
$ws.Range("M10").Value = 1.242987
$ws.Range("N10").Value = 2.485974
$ws.Range("O10").Value = 0.01287760861197995
$ws.Range("P10").Value = 0.008622083011467191
$ws.Range("Q10").Value = 5.624269234916
$ws.Range("R10").Value = 33.745615409496
$ws.Range("S10").Value = 0.002289752005157469
$ws.Range("T10").Value = 0.001640555137866101
$ws.Range("I11").Value = 0.1778087899819636
$ws.Range("J11").Value = 0.1902736421911268
$ws.Range("M11").Value = 1.174528666666667
$ws.Range("N11").Value = 3.523586
$ws.Range("O11").Value = 0.01216836577766621
$ws.Range("P11").Value = 0.01222082410759068
$ws.Range("Q11").Value = 5.314508876971555
$ws.Range("R11").Value = 47.83057989274401
$ws.Range("S11").Value = 0.002163642394984765
$ws.Range("T11").Value = 0.002325300713528406
$ws.Range("G12").Value = 4.544410333333333
$ws.Range("H12").Value = 13.633231
$ws.Range("I12").Value = 0.1785793547661169
$ws.Range("J12").Value = 0.1910982255429393
$ws.Range("M12").Value = 86.066935
$ws.Range("N12").Value = 258.200805
$ws.Range("O12").Value = 0.8916716774694496
$ws.Range("P12").Value = 0.8955157110805073
$ws.Range("Q12").Value = 391.1234687723284
$ws.Range("R12").Value = 3520.111218950955
$ws.Range("S12").Value = 0.1592341528257154
$ws.Range("T12").Value = 0.1711314633333085
$ws.Range("G13").Value = 4.544410333333333
$ws.Range("H13").Value = 13.633231
$ws.Range("I13").Value = 0.1785793547661169
$ws.Range("J13").Value = 0.1910982255429393
$ws.Range("O13").Value = 0.003820894467605101
$ws.Range("P13").Value = 0.003837366502243974
$ws.Range("Q13").Value = 1.675999738181556
$ws.Range("R13").Value = 15.083997643634
$ws.Range("S13").Value = 0.0006823328686543446
$ws.Range("T13").Value = 0.0007333139293367392
$ws.Range("G14").Value = 4.544410333333333
$ws.Range("H14").Value = 13.633231
$ws.Range("I14").Value = 0.1785793547661169
$ws.Range("J14").Value = 0.1910982255429393
$ws.Range("M14").Value = 7.669867666666666
$ws.Range("N14").Value = 23.009603
$ws.Range("O14").Value = 0.07946145367329926
$ws.Range("P14").Value = 0.07980401529819077
$ws.Range("Q14").Value = 34.85502587969922
$ws.Range("R14").Value = 313.695232917293
$ws.Range("S14").Value = 0.01419017512575547
$ws.Range("T14").Value = 0.01525040571468584
$ws.Range("G15").Value = 4.544410333333333
$ws.Range("H15").Value = 13.633231
$ws.Range("I15").Value = 0.1785793547661169
$ws.Range("J15").Value = 0.1910982255429393
$ws.Range("M15").Value = 1.242987
$ws.Range("N15").Value = 2.485974
$ws.Range("O15").Value = 0.01287760861197995
$ws.Range("P15").Value = 0.008622083011467191
$ws.Range("Q15").Value = 5.648642966999
$ws.Range("R15").Value = 33.891857801994
$ws.Range("S15").Value = 0.002299675036857969
$ws.Range("T15").Value = 0.001647664763975303
$ws.Range("G16").Value = 4.544410333333333
$ws.Range("H16").Value = 13.633231
$ws.Range("I16").Value = 0.1785793547661169
$ws.Range("J16").Value = 0.1910982255429393
$ws.Range("M16").Value = 1.174528666666667
$ws.Range("N16").Value = 3.523586
$ws.Range("O16").Value = 0.01216836577766621
$ws.Range("P16").Value = 0.01222082410759068
$ws.Range("Q16").Value = 5.337540209596222
$ws.Range("R16").Value = 48.037861886366
$ws.Range("S16").Value = 0.002173018909133731
$ws.Range("T16").Value = 0.002335377801632954
$ws.Range("G17").Value = 5.001220999999999
$ws.Range("H17").Value = 10.002442
$ws.Range("I17").Value = 0.1965304085046502
$ws.Range("J17").Value = 0.1402051294587592
$ws.Range("M17").Value = 86.066935
$ws.Range("N17").Value = 258.200805
$ws.Range("O17").Value = 0.8916716774694496
$ws.Range("P17").Value = 0.8955157110805073
$ws.Range("Q17").Value = 430.439762727635
$ws.Range("R17").Value = 2582.63857636581
$ws.Range("S17").Value = 0.1752405990250976
$ws.Range("T17").Value = 0.1255558962043953
$ws.Range("G18").Value = 5.001220999999999
$ws.Range("H18").Value = 10.002442
$ws.Range("I18").Value = 0.1965304085046502
$ws.Range("J18").Value = 0.1402051294587592
$ws.Range("O18").Value = 0.003820894467605101
$ws.Range("P18").Value = 0.003837366502243974
$ws.Range("Q18").Value = 1.844473643831333
$ws.Range("R18").Value = 11.066841862988
$ws.Range("S18").Value = 0.0007509219505715884
$ws.Range("T18").Value = 0.0005380184672278222
$ws.Range("G19").Value = 5.001220999999999
$ws.Range("H19").Value = 10.002442
$ws.Range("I19").Value = 0.1965304085046502
$ws.Range("J19").Value = 0.1402051294587592
$ws.Range("M19").Value = 7.669867666666666
$ws.Range("N19").Value = 23.009603
$ws.Range("O19").Value = 0.07946145367329926
$ws.Range("P19").Value = 0.07980401529819077
$ws.Range("Q19").Value = 38.35870324175433
$ws.Range("R19").Value = 230.152219450526
$ws.Range("S19").Value = 0.01561659195078684
$ws.Range("T19").Value = 0.01118893229621164
$ws.Range("G20").Value = 5.001220999999999
$ws.Range("H20").Value = 10.002442
$ws.Range("I20").Value = 0.1965304085046502
$ws.Range("J20").Value = 0.1402051294587592
$ws.Range("M20").Value = 1.242987
$ws.Range("N20").Value = 2.485974
$ws.Range("O20").Value = 0.01287760861197995
$ws.Range("P20").Value = 0.008622083011467191
$ws.Range("Q20").Value = 6.216452687126998
$ws.Range("R20").Value = 24.86581074850799
$ws.Range("S20").Value = 0.002530841681075421
$ws.Range("T20").Value = 0.001208860264826926
$ws.Range("G21").Value = 5.001220999999999
$ws.Range("H21").Value = 10.002442
$ws.Range("I21").Value = 0.1965304085046502
$ws.Range("J21").Value = 0.1402051294587592
$ws.Range("M21").Value = 1.174528666666667
$ws.Range("N21").Value = 3.523586
$ws.Range("O21").Value = 0.01216836577766621
$ws.Range("P21").Value = 0.01222082410759068
$ws.Range("Q21").Value = 5.874077432835332
$ws.Range("R21").Value = 35.24446459701199
$ws.Range("S21").Value = 0.002391453897118746
$ws.Range("T21").Value = 0.001713422226097476
$ws.Range("G22").Value = 8.684844333333333
$ws.Range("H22").Value = 26.054533
$ws.Range("I22").Value = 0.341283859407392
$ws.Range("J22").Value = 0.3652087332525911
$ws.Range("M22").Value = 86.066935
$ws.Range("N22").Value = 258.200805
$ws.Range("O22").Value = 0.8916716774694496
$ws.Range("P22").Value = 0.8955157110805073
$ws.Range("Q22").Value = 747.4779327221182
$ws.Range("R22").Value = 6727.301394499065
$ws.Range("S22").Value = 0.304313151411037
$ws.Range("T22").Value = 0.3270501584515055
$ws.Range("G23").Value = 8.684844333333333
$ws.Range("H23").Value = 26.054533
$ws.Range("I23").Value = 0.341283859407392
$ws.Range("J23").Value = 0.3652087332525911
$ws.Range("O23").Value = 0.003820894467605101
$ws.Range("P23").Value = 0.003837366502243974
$ws.Range("Q23").Value = 3.203011119406888
$ws.Range("R23").Value = 28.827100074662
$ws.Range("S23").Value = 0.001304009610292621
$ws.Range("T23").Value = 0.001401439759310448
$ws.Range("G24").Value = 8.684844333333333
$ws.Range("H24").Value = 26.054533
$ws.Range("I24").Value = 0.341283859407392
$ws.Range("J24").Value = 0.3652087332525911
$ws.Range("M24").Value = 7.669867666666666
$ws.Range("N24").Value = 23.009603
$ws.Range("O24").Value = 0.07946145367329926
$ws.Range("P24").Value = 0.07980401529819077
$ws.Range("Q24").Value = 66.61160674226655
$ws.Range("R24").Value = 599.5044606803989
$ws.Range("S24").Value = 0.02711891158374526
$ws.Range("T24").Value = 0.02914512333552265
$ws.Range("G25").Value = 8.684844333333333
$ws.Range("H25").Value = 26.054533
$ws.Range("I25").Value = 0.341283859407392
$ws.Range("J25").Value = 0.3652087332525911
$ws.Range("M25").Value = 1.242987
$ws.Range("N25").Value = 2.485974
$ws.Range("O25").Value = 0.01287760861197995
$ws.Range("P25").Value = 0.008622083011467191
$ws.Range("Q25").Value = 10.795148603357
$ws.Range("R25").Value = 64.770891620142
$ws.Range("S25").Value = 0.004394919967034385
$ws.Range("T25").Value = 0.003148860014616619
$ws.Range("G26").Value = 8.684844333333333
$ws.Range("H26").Value = 26.054533
$ws.Range("I26").Value = 0.341283859407392
$ws.Range("J26").Value = 0.3652087332525911
$ws.Range("M26").Value = 1.174528666666667
$ws.Range("N26").Value = 3.523586
$ws.Range("O26").Value = 0.01216836577766621
$ws.Range("P26").Value = 0.01222082410759068
$ws.Range("Q26").Value = 10.20059863503755
$ws.Range("R26").Value = 91.805387715338
$ws.Range("S26").Value = 0.004152866835282756
$ws.Range("T26").Value = 0.00446315169163592
